$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.235.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.603.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.96%  "

# Row 6
$ws.Range("E6").Value = "  +1.16%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("E9").Value = "  -1.60%  "

# Row 10
$ws.Range("E10").Value = "  +1.98%  "

# Row 11
$ws.Range("E11").Value = "  +1.31%  "

# Row 12
$ws.Range("E12").Value = "  +0.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.060.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.162.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.52%  "

# Row 15
$ws.Range("E15").Value = "  +0.43%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.601.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("E20").Value = "  -0.72%  "

# Row 21
$ws.Range("E21").Value = "  -1.86%  "

# Row 22
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.35%  "

# Row 24
$ws.Range("E24").Value = "  +0.89%  "

# Row 25
$ws.Range("E25").Value = "  -1.53%  "

# Row 26
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0753"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.96%  "

# Row 29
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("E30").Value = "  +7.66%  "

# Row 31
$ws.Range("E31").Value = "  -2.03%  "

# Row 32
$ws.Range("E32").Value = "  -0.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "

# Row 34
$ws.Range("E34").Value = "  +0.28%  "

# Row 35
$ws.Range("E35").Value = "  -0.86%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.25%  "

# Row 37
$ws.Range("E37").Value = "  -0.25%  "

# Row 38
$ws.Range("E38").Value = "  +0.60%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.49%  "

# Row 40
$ws.Range("E40").Value = "  +1.59%  "

# Row 41
$ws.Range("E41").Value = "  +0.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "274.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("E43").Value = "  +1.40%  "

# Row 44
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("E45").Value = "  +0.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0523"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47
$ws.Range("E47").Value = "  -1.84%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.33%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.31%  "

# Row 50
$ws.Range("E50").Value = "  +0.36%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
